$wb = $excel.ActiveWorkbook

# Rename variable labels (case fix) across every sheet in the workbook.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("dissolved_oxygen_uncorrected_mg_per_L", "dissolved_oxygen_uncorrected_mg_per_l")
    $ws.Cells.Replace("temperature_degree_C", "temperature_degree_c")
}

# Data correction on the "spike" sheet: spike_high threshold for
# temperature_degree_c changes from 5 to 4.
$spike = $wb.Worksheets.Item("spike")
$spike.Range("B6").Value = 4

# Restore per-sheet selections to match the saved workbook state.
$climatology = $wb.Worksheets.Item("climatology")
$climatology.Activate()
$climatology.Range("A14:A17").Select()

$grossrange = $wb.Worksheets.Item("grossrange")
$grossrange.Activate()
$grossrange.Range("B17").Select()

$seasons = $wb.Worksheets.Item("seasons")
$seasons.Activate()
$seasons.Range("F4").Select()

$spike.Activate()
$spike.Range("B7").Select()
